$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.099.83'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -3.04%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.647.51'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  -4.90%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9983'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.44'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -4.79%  '

$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4748'
$ws.Range("D7").Style = "Normal"

$ws.Range("E7").Value = '  -4.99%  '

$ws.Range("E8").Value = '  -5.47%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06100'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -0.80%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07035'
$ws.Range("D10").Style = "Normal"

$ws.Range("E10").Value = '  -2.77%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.640.34'
$ws.Range("D11").Style = "Normal"

$ws.Range("E11").Value = '  -5.67%  '

$ws.Range("E12").Value = '  -4.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.303'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  -9.53%  '

$ws.Range("E14").Value = '  -11.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '73.49'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  -4.46%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9993'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -0.19%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9992'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -0.19%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.096.26'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  -3.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.31'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = '  -4.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006609'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  -2.75%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.850.02'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -5.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.311'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  -6.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.437'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = '  -3.55%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.251'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = '  -3.52%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.00'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").Value = '  +0.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.92'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -1.93%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.377'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -2.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '103.36'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  -1.71%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.632'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -8.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.885'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  -1.80%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07587'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -6.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.524'
$ws.Range("D32").Style = "Normal"

$ws.Range("E32").Value = '  -4.27%  '

$ws.Range("E33").Value = '  -0.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04246'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  -10.25%  '

$ws.Range("E35").Value = '  -3.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9347'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -5.84%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5926'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -2.32%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.583'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  -5.47%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8555'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9992'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").Value = '  -0.19%  '

$ws.Range("E41").Value = '  -7.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '98.42'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  -1.47%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.782'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = '  -7.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3676'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -5.74%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.618'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -7.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1098'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").Value = '  -6.43%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05218'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -0.96%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.069'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  -3.82%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '28.97'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  -5.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9991'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  -0.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  +0.05%  '
